# Build site at 2023-04-12 14:53:07 UTC
# Update LOT2013 syllabus sheet: fill in the previously-missing Portuguese
# "Objetivos" / "Programa resumido" / "Programa" / "Bibliografia" content,
# fix the objectives text that had been mistakenly duplicated from
# "Walter de Carvalho", and append the new bibliography row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# New / corrected text blocks
# ---------------------------------------------------------------------
$objetivosPt   = "Capacitar os alunos para a aplicação de conhecimentos da engenharia química na solução de problemas que se apresentam na implantação e otimização de processos biotecnológicos, com ênfase em: 1) esterilização de equipamentos, meios e ar e; 2) cinética e estequiometria do crescimento microbiano e da formação de produtos."
$walter        = "5840876 - Walter de Carvalho"
$rotProgResum  = "Programa resumido:"
$progResumPt   = "Processo biotecnológico genérico, esterilização de equipamentos, esterilização de meios por aquecimento com vapor, esterilização de ar por filtração, cinética e estequiometria do crescimento microbiano e da formação de produtos."
$rotShortSyl   = "Short syllabus:"
$shortSylEn    = "Generic biotechnological process, equipment sterilization, media sterilization by heating, air sterilization by filtration, kinetics and stoichiometry of microbial growth and products formation."
$rotPrograma   = "Programa:"
$programaPt    = "1. Processo biotecnológico genérico: representação esquemática; descrição das etapas principais.2. Esterilização de equipamentos: terminologia; esterilização por agentes físicos; esterilização por agentes químicos.3. Esterilização de meios por aquecimento com vapor: cinética da destruição térmica de microrganismos; destruição de nutrientes do meio; cálculo do tempo de esterilização por processo descontínuo; cálculo do tempo de esterilização por processo contínuo.4. Esterilização de ar por filtração: aerossóis microbianos; amostradores; dimensionamento de filtros fibrosos; dimensionamento de filtros de membranas.5. Cinética e estequiometria do crescimento microbiano e da formação de produtos: velocidades de transformação e fatores de conversão; classificação dos processos fermentativos em função das velocidades de crescimento celular e formação de produtos; influência da concentração do substrato sobre a velocidade de crescimento celular; estequiometria do crescimento microbiano e da formação de produtos."
$rotSyllabus   = "Syllabus:"
$syllabusEn    = "1.Generic biotechnological process: schematic representation; description of the main phases.2.Equipment sterilization: terminology, sterilization by physical agents, sterilization by chemical agents.3.Media sterilization by heating: kinetics of thermic destruction of microorganisms; destruction of medium nutrients; calculation of sterilization time by discontinuous process; calculation of sterilization time by continuous process.4.Air sterilization by filtration: microbial aerosols; air samplers; dimensioning of fibrous filters; dimensioning of membrane filters.5.Kinetics and stoichiometry of microbial growth and products formation: definition of velocities of transformation and conversion factors; classification of fermentations as a function of the velocities of cell growth and products formation; influence of substrate concentration on cell growth; stoichiometry of microbial growth and products formation."
$rotAvaliacao  = "Avaliação:"
$rotMetodo     = "Método:"
$metodoTxt     = "Os alunos serão avaliados formalmente por duas provas escritas (P1 e P2), sendo a segunda prova (P2) com peso 2."
$rotCriterio   = "Critério:"
$criterioTxt   = "A nota final (NF) será calculada como: NF=(P1+(P2×2))/3. Serão aprovados os alunos que obtiverem NF maior ou igual 5,0."
$rotNorma      = "Norma de recuperação:"
$normaTxt      = "Será oferecido um programa de recuperação avaliado por uma prova escrita final (PR).`nA média de recuperação (MR) será calculada como: MR=(NF+PR)/2. Serão aprovados os alunos que obtiverem MR maior ou igual a 5,0."
$rotBiblio     = "Bibliografia:"
$biblioTxt     = "BAILEY, J.E., OLLIS D.F. Biochemical Engineering Fundamentals. 2nd edition, New York: McGraw Hill, 1986. ISBN: 978-0070032125.BORZANI, W., SCHMIDELL, W., LIMA, U.A., AQUARONE, E. Biotecnologia Industrial Fundamentos (Vol 1). São Paulo: Edgard Blucher Ltda, 2001.m ISBN: 9788521202783.DORAN P.M. Bioprocess Engineering Principles, 1st edition, San Diego: Academic Press, 1995. ISBN: 978-0080528120.KATOH, S., HORIUCHI, J., YOSHIDA, F. Biochemical Engineering: A Textbook for Engineers, Chemists and Biologists, 2nd, Completely Revised and Enlarged Edition. Weinheim/Germany: Wiley-VCH, 2015. ISBN: 978-3527338047.SCHMIDELL, W., LIMA, U.A., AQUARONE, E., BORZANI, W. Biotecnologia Industrial Engenharia Bioquímica (Vol 2), São Paulo: Edgard Blucher Ltda, 2001. ISBN: 9788521202790."
$rotRequisitos = "Requisitos:"
$lot2028Txt    = "LOT2028 -  Tecnologia de Processos Fermentativos  (Requisito fraco)`n"

# ---------------------------------------------------------------------
# Row 10: fix the Objetivos/Objectives body text (was wrongly showing
# the "Walter de Carvalho" string that belongs to "Docentes responsáveis").
# ---------------------------------------------------------------------
$ws.Range("B10").Value2 = $objetivosPt
$ws.Range("C10").Value2 = $objetivosPt

# ---------------------------------------------------------------------
# Row 13 used to hold "Programa resumido:" / "Semestral". It becomes the
# (new) home of the "Docentes responsáveis" value and loses its own
# custom row height (back to sheet default) and its A-label.
# ---------------------------------------------------------------------
$ws.Range("A13").Clear()
$ws.Range("B13").Value2 = $walter
$ws.Range("C13").Value2 = $walter
$ws.Rows.Item(13).AutoFit()

# ---------------------------------------------------------------------
# Rows 14-22: every row shifts down one slot from where it used to be,
# while row 14 gains fresh "Programa resumido" content and row 16 / 22
# gain fresh "Programa" / "Bibliografia" bodies.
# ---------------------------------------------------------------------
$ws.Range("A14").Value2 = $rotProgResum
$ws.Range("B14").Value2 = $progResumPt
$ws.Range("C14").Value2 = $progResumPt
$ws.Rows.Item(14).RowHeight = 60

$ws.Range("A15").Value2 = $rotShortSyl
$ws.Range("B15").Value2 = $shortSylEn
$ws.Range("C15").Value2 = $shortSylEn
$ws.Rows.Item(15).RowHeight = 60

$ws.Range("A16").Value2 = $rotPrograma
$ws.Range("B16").Value2 = $programaPt
$ws.Range("C16").Value2 = $programaPt
$ws.Rows.Item(16).RowHeight = 120

$ws.Range("A17").Value2 = $rotSyllabus
$ws.Range("B17").Value2 = $syllabusEn
$ws.Range("C17").Value2 = $syllabusEn
$ws.Rows.Item(17).RowHeight = 120

$ws.Range("A18").Value2 = $rotAvaliacao
$ws.Rows.Item(18).AutoFit()

$ws.Range("A19").Value2 = $rotMetodo
$ws.Range("B19").Value2 = $metodoTxt
$ws.Range("C19").Value2 = $metodoTxt
$ws.Rows.Item(19).RowHeight = 60

$ws.Range("A20").Value2 = $rotCriterio
$ws.Range("B20").Value2 = $criterioTxt
$ws.Range("C20").Value2 = $criterioTxt
$ws.Rows.Item(20).RowHeight = 60

$ws.Range("A21").Value2 = $rotNorma
$ws.Range("B21").Value2 = $normaTxt
$ws.Range("C21").Value2 = $normaTxt
$ws.Rows.Item(21).RowHeight = 60

$ws.Range("A22").Value2 = $rotBiblio
$ws.Range("B22").Value2 = $biblioTxt
$ws.Range("C22").Value2 = $biblioTxt
$ws.Rows.Item(22).RowHeight = 120

# ---------------------------------------------------------------------
# Row 23 used to carry the "LOT2028..." requisito text in B/C; that now
# moves to the brand-new row 24, and row 23 becomes the "Requisitos:"
# label row (style copied from the A column above it).
# ---------------------------------------------------------------------
$ws.Range("A22").Copy() | Out-Null
$ws.Range("A23").PasteSpecial(-4122) | Out-Null
$ws.Range("A23").Value2 = $rotRequisitos
$ws.Rows.Item(23).AutoFit()

$ws.Range("B23").Copy() | Out-Null
$ws.Range("B24").PasteSpecial(-4122) | Out-Null
$ws.Range("C23").Copy() | Out-Null
$ws.Range("C24").PasteSpecial(-4122) | Out-Null
$ws.Range("B24").Value2 = $lot2028Txt
$ws.Range("C24").Value2 = $lot2028Txt
$ws.Rows.Item(24).RowHeight = 30

$ws.Application.CutCopyMode = $false

# ---------------------------------------------------------------------
# Columns: column A (width 30.7109375) needs to stand on its own instead
# of sharing a single <col> entry with column B (width 60.7109375).
# ---------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = $ws.Columns.Item(2).ColumnWidth
